$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 158-159; this pushes the existing rows
# 158..179 down to 160..181 and keeps their formatting/values intact.
$ws.Rows("158:159").Insert()

# Row 158 - new record (week of 2021-11-04), Española / Primera
$ws.Cells.Item(158, 1).Value = 2
$ws.Cells.Item(158, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(158, 3).Value = "Coquimbo"
$ws.Cells.Item(158, 4).Value = 44504
$ws.Cells.Item(158, 5).Value = 4
$ws.Cells.Item(158, 6).Value = 100112013
$ws.Cells.Item(158, 7).Value = "Alcachofa"
$ws.Cells.Item(158, 8).Value = "Española"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 1000
$ws.Cells.Item(158, 11).Value = 4500
$ws.Cells.Item(158, 12).Value = 5000
$ws.Cells.Item(158, 13).Value = 4750
$ws.Cells.Item(158, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(158, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(158, 16).Value = 158
$ws.Cells.Item(158, 17).Value = 30
$ws.Cells.Item(158, 18).Value = "Hortaliza"

# Row 159 - new record (week of 2021-11-04), Madrigal / Primera
$ws.Cells.Item(159, 1).Value = 2
$ws.Cells.Item(159, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(159, 3).Value = "Coquimbo"
$ws.Cells.Item(159, 4).Value = 44504
$ws.Cells.Item(159, 5).Value = 4
$ws.Cells.Item(159, 6).Value = 100112013
$ws.Cells.Item(159, 7).Value = "Alcachofa"
$ws.Cells.Item(159, 8).Value = "Madrigal"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 1100
$ws.Cells.Item(159, 11).Value = 4500
$ws.Cells.Item(159, 12).Value = 5000
$ws.Cells.Item(159, 13).Value = 4750
$ws.Cells.Item(159, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(159, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(159, 16).Value = 119
$ws.Cells.Item(159, 17).Value = 40
$ws.Cells.Item(159, 18).Value = "Hortaliza"
